# Update the price list date and the four support prices on Hoja1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds a date serial number (formatted as a date). Move it forward
# from 45406 (2024-04-24) to 45436 (2024-05-24).
$ws.Range("A1").Value = 45436

# Update the four unit prices in column D (rows 23-26) to 844.
$ws.Range("D23").Value = 844
$ws.Range("D24").Value = 844
$ws.Range("D25").Value = 844
$ws.Range("D26").Value = 844
